# Apply "update to try logic" changes to the point table (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - Team C: Won 3->4, Drawn 1->0, Goal Scored 1->2
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 0
$ws.Range("F3").Value = 2

# Row 4 - Team D: Goal Scored 0->5, Penalty Points 0->-6
$ws.Range("F4").Value = 5
$ws.Range("H4").Value = -6

# Row 5 - Team E: Won 1->2, Drawn 3->2, Goal Concedered 5->2, Penalty Points 0->-3
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = -3

# Update the active selection to J6 as recorded in the workbook view
$ws.Activate()
$ws.Range("J6").Select()
